$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$updates = @(
    @{Row=17; Col=8; Value=1371484.2}
    @{Row=17; Col=10; Value=1395516.2}
    @{Row=17; Col=12; Value=4186548.6}
    @{Row=17; Col=14; Value=-4186884.6}
    @{Row=19; Col=8; Value=437.33334}
    @{Row=19; Col=9; Value=225}
    @{Row=19; Col=10; Value=498}
    @{Row=19; Col=11; Value=225}
    @{Row=19; Col=12; Value=498}
    @{Row=19; Col=13; Value=-50}
    @{Row=19; Col=14; Value=-848}
    @{Row=40; Col=8; Value=4369.72}
    @{Row=40; Col=9; Value=1818.6666}
    @{Row=40; Col=11; Value=1818.6666}
    @{Row=40; Col=13; Value=-1643.6666}
    @{Row=62; Col=8; Value=4213.077}
    @{Row=62; Col=9; Value=4057.5}
    @{Row=62; Col=11; Value=4057.5}
    @{Row=62; Col=13; Value=-3433.5}
    @{Row=65; Col=8; Value=4213.077}
    @{Row=65; Col=9; Value=4057.5}
    @{Row=65; Col=11; Value=20287.5}
    @{Row=65; Col=13; Value=-17167.5}
    @{Row=112; Col=8; Value=502159.5}
    @{Row=112; Col=9; Value=724.75}
    @{Row=112; Col=11; Value=2174.25}
    @{Row=112; Col=13; Value=-1066.25}
    @{Row=135; Col=8; Value=1749.0526}
    @{Row=135; Col=9; Value=1719.5883}
    @{Row=135; Col=11; Value=15476.2947}
    @{Row=135; Col=13; Value=-12941.2947}
    @{Row=137; Col=8; Value=2510.6177}
    @{Row=137; Col=9; Value=2398.5417}
    @{Row=137; Col=11; Value=7195.625100000001}
    @{Row=137; Col=13; Value=-4645.625100000001}
    @{Row=141; Col=8; Value=5149}
    @{Row=141; Col=9; Value=3950}
    @{Row=141; Col=10; Value=5748.5}
    @{Row=141; Col=11; Value=11850}
    @{Row=141; Col=12; Value=17245.5}
    @{Row=141; Col=13; Value=-6670}
    @{Row=141; Col=14; Value=-27605.5}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item("ARM")
$updates = @(
    @{Row=2; Col=8; Value=1842.9333}
    @{Row=2; Col=9; Value=1884.9166}
    @{Row=2; Col=11; Value=1884.9166}
    @{Row=2; Col=13; Value=-1771.9166}
    @{Row=32; Col=8; Value=7926.0483}
    @{Row=32; Col=9; Value=8163.271}
    @{Row=32; Col=11; Value=8163.271}
    @{Row=32; Col=13; Value=-7876.271}
    @{Row=45; Col=8; Value=4546.778}
    @{Row=45; Col=9; Value=3986}
    @{Row=45; Col=11; Value=3986}
    @{Row=45; Col=13; Value=-3609}
    @{Row=116; Col=8; Value=1842.9333}
    @{Row=116; Col=9; Value=1884.9166}
    @{Row=116; Col=11; Value=1884.9166}
    @{Row=116; Col=13; Value=409.0834}
    @{Row=132; Col=8; Value=2138504.2}
    @{Row=132; Col=9; Value=2483008.2}
    @{Row=132; Col=11; Value=7449024.600000001}
    @{Row=132; Col=13; Value=-7446494.600000001}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item("BSM")
$updates = @(
    @{Row=3; Col=8; Value=1842.9333}
    @{Row=3; Col=9; Value=1884.9166}
    @{Row=3; Col=11; Value=1884.9166}
    @{Row=3; Col=13; Value=-1770.9166}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item("CRP")
$updates = @(
    @{Row=4; Col=8; Value=111117780}
    @{Row=4; Col=9; Value=5000}
    @{Row=4; Col=10; Value=125006870}
    @{Row=4; Col=11; Value=5000}
    @{Row=4; Col=12; Value=125006870}
    @{Row=4; Col=13; Value=-4888}
    @{Row=4; Col=14; Value=-125007094}
    @{Row=81; Col=8; Value=79666.664}
    @{Row=81; Col=10; Value=79666.664}
    @{Row=81; Col=12; Value=79666.664}
    @{Row=81; Col=14; Value=-81662.664}
    @{Row=84; Col=8; Value=79666.664}
    @{Row=84; Col=10; Value=79666.664}
    @{Row=84; Col=12; Value=238999.992}
    @{Row=84; Col=14; Value=-248983.992}
    @{Row=122; Col=8; Value=2912.9092}
    @{Row=122; Col=9; Value=2959.84}
    @{Row=122; Col=10; Value=2766.25}
    @{Row=122; Col=11; Value=8879.52}
    @{Row=122; Col=12; Value=8298.75}
    @{Row=122; Col=13; Value=-6429.52}
    @{Row=122; Col=14; Value=-13198.75}
    @{Row=132; Col=8; Value=1179757.6}
    @{Row=132; Col=9; Value=2108290.2}
    @{Row=132; Col=11; Value=6324870.600000001}
    @{Row=132; Col=13; Value=-6322340.600000001}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item("CUL")
$updates = @(
    @{Row=4; Col=8; Value=12731809}
    @{Row=4; Col=9; Value=14575097}
    @{Row=4; Col=10; Value=9736466}
    @{Row=4; Col=11; Value=43725291}
    @{Row=4; Col=12; Value=29209398}
    @{Row=4; Col=13; Value=-43725179}
    @{Row=4; Col=14; Value=-29209622}
    @{Row=45; Col=8; Value=1000}
    @{Row=45; Col=10; Value=1000}
    @{Row=45; Col=12; Value=3000}
    @{Row=45; Col=14; Value=-4064}
    @{Row=63; Col=8; Value=5956.9165}
    @{Row=63; Col=9; Value=3329.3333}
    @{Row=63; Col=10; Value=6832.778}
    @{Row=63; Col=11; Value=9987.999899999999}
    @{Row=63; Col=12; Value=20498.334}
    @{Row=63; Col=13; Value=-9238.999899999999}
    @{Row=63; Col=14; Value=-21996.334}
    @{Row=66; Col=8; Value=5956.9165}
    @{Row=66; Col=9; Value=3329.3333}
    @{Row=66; Col=10; Value=6832.778}
    @{Row=66; Col=11; Value=29963.9997}
    @{Row=66; Col=12; Value=61495.002}
    @{Row=66; Col=13; Value=-26219.9997}
    @{Row=66; Col=14; Value=-68983.00200000001}
    @{Row=98; Col=8; Value=2188.7778}
    @{Row=98; Col=9; Value=699.6667}
    @{Row=98; Col=11; Value=2099.0001}
    @{Row=98; Col=13; Value=-601.0001000000002}
    @{Row=105; Col=8; Value=14571}
    @{Row=105; Col=10; Value=14571}
    @{Row=105; Col=12; Value=43713}
    @{Row=105; Col=14; Value=-48955}
    @{Row=109; Col=8; Value=4269.6875}
    @{Row=109; Col=9; Value=1258.4445}
    @{Row=109; Col=11; Value=3775.3335}
    @{Row=109; Col=13; Value=-2735.3335}
    @{Row=121; Col=8; Value=5841.2}
    @{Row=121; Col=10; Value=5841.2}
    @{Row=121; Col=12; Value=17523.6}
    @{Row=121; Col=14; Value=-20143.6}
    @{Row=129; Col=8; Value=1308.4615}
    @{Row=129; Col=9; Value=753.5}
    @{Row=129; Col=11; Value=2260.5}
    @{Row=129; Col=13; Value=2739.5}
    @{Row=131; Col=8; Value=1598.6875}
    @{Row=131; Col=10; Value=1601.5}
    @{Row=131; Col=12; Value=4804.5}
    @{Row=131; Col=14; Value=-14884.5}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item("GSM")
$updates = @(
    @{Row=97; Col=8; Value=1471.3422}
    @{Row=97; Col=9; Value=541.9048}
    @{Row=97; Col=10; Value=2619.4707}
    @{Row=97; Col=11; Value=541.9048}
    @{Row=97; Col=12; Value=2619.4707}
    @{Row=97; Col=13; Value=-45.90480000000002}
    @{Row=97; Col=14; Value=-3611.4707}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item("LTW")
$updates = @(
    @{Row=6; Col=8; Value=77499.5}
    @{Row=6; Col=10; Value=77499.5}
    @{Row=6; Col=12; Value=77499.5}
    @{Row=6; Col=14; Value=-77723.5}
    @{Row=7; Col=8; Value=3791.0417}
    @{Row=7; Col=9; Value=3410.0527}
    @{Row=7; Col=10; Value=5238.8}
    @{Row=7; Col=11; Value=3410.0527}
    @{Row=7; Col=12; Value=5238.8}
    @{Row=7; Col=13; Value=-3298.0527}
    @{Row=7; Col=14; Value=-5462.8}
    @{Row=22; Col=8; Value=2175.5652}
    @{Row=22; Col=9; Value=721.1111}
    @{Row=22; Col=10; Value=3110.5715}
    @{Row=22; Col=11; Value=721.1111}
    @{Row=22; Col=12; Value=3110.5715}
    @{Row=22; Col=13; Value=-426.1111}
    @{Row=22; Col=14; Value=-3700.5715}
    @{Row=27; Col=8; Value=2175.5652}
    @{Row=27; Col=9; Value=721.1111}
    @{Row=27; Col=10; Value=3110.5715}
    @{Row=27; Col=11; Value=721.1111}
    @{Row=27; Col=12; Value=3110.5715}
    @{Row=27; Col=13; Value=-614.1111}
    @{Row=27; Col=14; Value=-3324.5715}
    @{Row=68; Col=8; Value=2361.125}
    @{Row=68; Col=9; Value=2478.6}
    @{Row=68; Col=10; Value=2165.3333}
    @{Row=68; Col=11; Value=2478.6}
    @{Row=68; Col=12; Value=2165.3333}
    @{Row=68; Col=13; Value=-1729.6}
    @{Row=68; Col=14; Value=-3663.3333}
    @{Row=71; Col=8; Value=2361.125}
    @{Row=71; Col=9; Value=2478.6}
    @{Row=71; Col=10; Value=2165.3333}
    @{Row=71; Col=11; Value=12393}
    @{Row=71; Col=12; Value=10826.6665}
    @{Row=71; Col=13; Value=-8649}
    @{Row=71; Col=14; Value=-18314.6665}
    @{Row=92; Col=8; Value=98888}
    @{Row=92; Col=10; Value=98888}
    @{Row=92; Col=12; Value=98888}
    @{Row=92; Col=14; Value=-103880}
    @{Row=126; Col=8; Value=3791.0417}
    @{Row=126; Col=9; Value=3410.0527}
    @{Row=126; Col=10; Value=5238.8}
    @{Row=126; Col=11; Value=10230.1581}
    @{Row=126; Col=12; Value=15716.4}
    @{Row=126; Col=13; Value=-7760.158100000001}
    @{Row=126; Col=14; Value=-20656.4}
    @{Row=132; Col=8; Value=2915.766}
    @{Row=132; Col=9; Value=2873.228}
    @{Row=132; Col=11; Value=8619.684000000001}
    @{Row=132; Col=13; Value=-6089.684000000001}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

$ws = $wb.Worksheets.Item("WVR")
$updates = @(
    @{Row=81; Col=8; Value=48653.043}
    @{Row=81; Col=10; Value=8061.364}
    @{Row=81; Col=12; Value=16122.728}
    @{Row=81; Col=14; Value=-18244.728}
    @{Row=84; Col=8; Value=48653.043}
    @{Row=84; Col=10; Value=8061.364}
    @{Row=84; Col=12; Value=80613.64}
    @{Row=84; Col=14; Value=-91221.64}
    @{Row=107; Col=8; Value=1220.8462}
    @{Row=107; Col=9; Value=997.36365}
    @{Row=107; Col=10; Value=2450}
    @{Row=107; Col=11; Value=2992.09095}
    @{Row=107; Col=12; Value=7350}
    @{Row=107; Col=13; Value=-1072.09095}
    @{Row=107; Col=14; Value=-11190}
    @{Row=126; Col=8; Value=1672.12}
    @{Row=126; Col=9; Value=1218.579}
    @{Row=126; Col=11; Value=3655.737}
    @{Row=126; Col=13; Value=-1185.737}
)
foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}
